$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to text format before writing so that
# numeric-looking strings (e.g. "333.45") are kept as text, matching
# the source data which stores these as plain strings, not numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.764.05"
$ws.Range("E2").Value = "  +2.14%  "

$ws.Range("D3").Value = "2.110.71"
$ws.Range("E3").Value = "  +8.32%  "

$ws.Range("E4").Value = "  +0.32%  "

$ws.Range("D5").Value = "333.45"
$ws.Range("E5").Value = "  +3.97%  "

$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.31%  "

$ws.Range("D7").Value = "0.5284"
$ws.Range("E7").Value = "  +3.85%  "

$ws.Range("D8").Value = "0.4384"
$ws.Range("E8").Value = "  +7.63%  "

$ws.Range("D9").Value = "0.09007"
$ws.Range("E9").Value = "  +7.07%  "

$ws.Range("D10").Value = "45.91"
$ws.Range("E10").Value = "  +8.38%  "

$ws.Range("D11").Value = "1.179"
$ws.Range("E11").Value = "  +4.92%  "

$ws.Range("E12").Value = "  +3.48%  "

$ws.Range("D13").Value = "2.109.55"
$ws.Range("E13").Value = "  +8.62%  "

$ws.Range("D14").Value = "6.751"
$ws.Range("E14").Value = "  +4.89%  "

$ws.Range("D15").Value = "7.800"
$ws.Range("E15").Value = "  +6.35%  "

$ws.Range("D16").Value = "97.40"
$ws.Range("E16").Value = "  +4.34%  "

$ws.Range("E17").Value = "  -0.04%  "

$ws.Range("D18").Value = "0.00001126"
$ws.Range("E18").Value = "  +2.41%  "

$ws.Range("D19").Value = "0.06670"
$ws.Range("E19").Value = "  +2.36%  "

$ws.Range("D20").Value = "19.13"
$ws.Range("E20").Value = "  +2.69%  "

$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  +0.09%  "

$ws.Range("D22").Value = "6.360"
$ws.Range("E22").Value = "  +5.96%  "

$ws.Range("D23").Value = "30.825.19"
$ws.Range("E23").Value = "  +2.29%  "

$ws.Range("D24").Value = "12.16"
$ws.Range("E24").Value = "  +6.66%  "

$ws.Range("D25").Value = "2.357.71"
$ws.Range("E25").Value = "  +9.12%  "

$ws.Range("D26").Value = "2.260"
$ws.Range("E26").Value = "  +2.98%  "

$ws.Range("D27").Value = "22.81"
$ws.Range("E27").Value = "  +1.75%  "

$ws.Range("E28").Value = "  +9.29%  "

$ws.Range("D29").Value = "162.43"
$ws.Range("E29").Value = "  -0.31%  "

$ws.Range("D30").Value = "132.94"
$ws.Range("E30").Value = "  +2.62%  "

$ws.Range("E31").Value = "  +2.82%  "

$ws.Range("D32").Value = "0.1074"
$ws.Range("E32").Value = "  +2.45%  "

$ws.Range("D33").Value = "6.223"
$ws.Range("E33").Value = "  +3.39%  "

$ws.Range("D34").Value = "4.055"
$ws.Range("E34").Value = "  +7.23%  "

$ws.Range("E35").Value = "  +20.60%  "

$ws.Range("E36").Value = "  +5.65%  "

$ws.Range("D37").Value = "5.539"
$ws.Range("E37").Value = "  +3.54%  "

$ws.Range("D38").Value = "0.06738"
$ws.Range("E38").Value = "  +3.77%  "

$ws.Range("D39").Value = "9.517"
$ws.Range("E39").Value = "  +8.28%  "

$ws.Range("E40").Value = "  +8.11%  "

$ws.Range("D41").Value = "0.2271"
$ws.Range("E41").Value = "  +5.01%  "

$ws.Range("D42").Value = "0.6832"
$ws.Range("E42").Value = "  +4.49%  "

$ws.Range("D43").Value = "1.242"
$ws.Range("E43").Value = "  +1.44%  "

$ws.Range("D44").Value = "0.6453"
$ws.Range("E44").Value = "  +5.71%  "

$ws.Range("B45").Value = "Frax"
$ws.Range("C45").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D45").Value = "1.001"
$ws.Range("E45").Value = "  +0.36%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "14.04"
$ws.Range("E46").Value = "  +4.82%  "

$ws.Range("D47").Value = "2.231"
$ws.Range("E47").Value = "  +2.03%  "

$ws.Range("D48").Value = "3.671"
$ws.Range("E48").Value = "  +1.18%  "

$ws.Range("E49").Value = "  +4.88%  "

$ws.Range("D50").Value = "82.42"
$ws.Range("E50").Value = "  +4.83%  "

$ws.Range("D51").Value = "119.52"
$ws.Range("E51").Value = "  -2.81%  "

# Restore the default cell style on the Price column so no stray
# text-format style is left applied to the cells (matches source formatting).
$ws.Range("D2:D51").Style = "Normal"
